# Edit: Added J-Link adapter and 1st cut of software
# Adds two new BOM rows (SFH11-PBPC-D05-ST-BK 10-pin plug, and
# 20021121-00010C4LF 10-pin header) below the existing parts list,
# leaves two blank spacer rows, and relocates the "Total" row from
# row 38 down to row 41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Capture the existing "Total" row (currently row 38) before we
#        overwrite it with new part rows. ---
$totalLabelFormula = $ws.Range("D38").Formula
$totalSumFormula    = $ws.Range("E38").Formula

# --- 2. Clear the old Total row cells (C38:E38); we are about to reuse
#        row 38 for new part data, and will rebuild the Total row at 41. ---
$ws.Range("C38:E38").ClearContents()
$ws.Range("C38:E38").ClearFormats()

# --- 3. New row 37: SFH11-PBPC-D05-ST-BK (10-pin plug, through hole) ---
$ws.Range("A37").Value2 = "SFH11-PBPC-D05-ST-BK"
$ws.Range("A37").Style = "Hyperlink"

$ws.Range("C37").Value2 = "10-pin plug - 0.1"""

$ws.Range("D37").Value2 = 1

$ws.Range("E37").Value2 = 0.66

$ws.Range("F37").Value2 = "Through Hole"
$ws.Range("F37").HorizontalAlignment = -4131   # xlLeft

# --- 4. New row 38: 20021121-00010C4LF (10-pin header, SMT) ---
$ws.Range("A38").Value2 = "20021121-00010C4LF"
$ws.Range("A38").Style = "Hyperlink"

$ws.Range("B38").Value2 = "Amphenol"

$ws.Range("C38").Value2 = "10-pin header - .005"""

$ws.Range("D38").Value2 = 1

$ws.Range("E38").Value2 = 0.84

$ws.Range("F38").Value2 = "SMT"
$ws.Range("F38").HorizontalAlignment = -4131   # xlLeft

# --- 5. Two blank spacer rows (39, 40) - column A keeps the same
#        "Hyperlink-less" part-number formatting used throughout column A. ---
$ws.Range("A36").Copy() | Out-Null
$ws.Range("A39").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A39").ClearContents()
$ws.Range("A39").Style = "Hyperlink"

$ws.Range("A40").PasteSpecial(-4122) | Out-Null
$ws.Range("A40").ClearContents()
$ws.Range("A40").Style = "Hyperlink"

# --- 6. Rebuild the Total row at 41 (was row 38). ---
$ws.Range("C41").Style = "Normal"
$ws.Range("C41").HorizontalAlignment = -4152    # xlRight

$ws.Range("D41").Formula = $totalLabelFormula
$ws.Range("D41").Font.Bold = $true
$ws.Range("D41").HorizontalAlignment = -4108    # xlCenter

$ws.Range("E41").Formula = $totalSumFormula
$ws.Range("E41").NumberFormat = $ws.Range("E2").NumberFormat
$ws.Range("E41").HorizontalAlignment = -4108    # xlCenter

# --- 7. Hyperlinks for the two new part numbers (added in row order so
#        relationship ids come out rId16 / rId17). ---
$ws.Hyperlinks.Add($ws.Range("A37"), "https://www.digikey.com/product-detail/en/sullins-connector-solutions/SFH11-PBPC-D05-ST-BK/S9194-ND/1990087", "", "", "SFH11-PBPC-D05-ST-BK") | Out-Null
$ws.Range("A37").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("A38"), "https://www.digikey.com/product-detail/en/amphenol-fci/20021121-00010C4LF/609-3695-1-ND/2209147", "", "", "20021121-00010C4LF") | Out-Null
$ws.Range("A38").Style = "Hyperlink"

# --- 8. Update the view: scrolled down, F38 selected. ---
$ws.Range("F38").Select() | Out-Null

Write-Output "NFC Toy BOM updated: added J-Link adapter parts + moved Total row to 41"
